$d = $word.ActiveDocument

function Replace-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $d.Range($p.Range.Start, $p.Range.End)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg)
}

# Paragraph 4: Campeonato (CampeonatoID, Nome, Local, DataInicio, DataFim, pais->Pais)
$inner4 = '<w:r><w:t>Campeonato (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Campeonato</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, Nome, Local, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>DataInicio</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>DataFim</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>pais-&gt;Pais</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>'
Replace-ParaXml 4 $inner4

# Paragraph 6: Categoria (CategoriaNome, NumMaxParticipantes, AlturaMinima, AlturaMaxima, PesoMinimo, PesoMaximo, Genero) //Triggers...
$inner6 = '<w:r><w:t>Categoria (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Categoria</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Nome</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>NumMaxParticipantes</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>AlturaMinima</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>AlturaMaxima</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>PesoMinimo</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>PesoMaximo</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Genero</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> //</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Triggers</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> -&gt; relat&#243;rio, necess&#225;rio para implementar a restri&#231;&#227;o altura e peso</w:t></w:r>'
Replace-ParaXml 6 $inner6

# Paragraph 9: Equipa (EquipaNome, Local)
$inner9 = '<w:r><w:t>Equipa (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Equipa</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Nome</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>, Local)</w:t></w:r>'
Replace-ParaXml 9 $inner9

# Paragraph 10: Fase (FaseNome)
$inner10 = '<w:r><w:t>Fase (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Fase</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Nome</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>'
Replace-ParaXml 10 $inner10

# Paragraph 11: Jurado (CC->Pessoa, juriID->Juri)  (juri/ID runs gain underline)
$inner11 = '<w:r><w:t>Jurado (</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CC</w:t></w:r>' + `
    '<w:r><w:t>-&gt;Pessoa</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>juri</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>-&gt;</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Juri</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>'
Replace-ParaXml 11 $inner11

# Paragraph 12: Juri (JuriID)
$inner12 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Juri</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Juri</w:t></w:r>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>'
Replace-ParaXml 12 $inner12

# Paragraphs 14 (Patrocinio) and 15 (Pessoa) need to be replaced together since the
# _GoBack bookmark moves from the Pessoa paragraph into the Patrocinio paragraph.
$p14 = $d.Paragraphs(14)
$p15 = $d.Paragraphs(15)
$full1415 = $d.Range($p14.Range.Start, $p15.Range.End)
$inner14 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Patrocinio</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
           '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
           '<w:proofErr w:type="spellStart"/>' + `
           '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Patrocinio</w:t></w:r>' + `
           '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Nome</w:t></w:r>' + `
           '<w:bookmarkEnd w:id="0"/>' + `
           '<w:proofErr w:type="spellEnd"/>' + `
           '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
           '<w:r><w:t>Valor)</w:t></w:r>'
$inner15 = '<w:r><w:t>Pessoa (</w:t></w:r>' + `
           '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>CC</w:t></w:r>' + `
           '<w:r><w:t xml:space="preserve">, Nome, Morada, </w:t></w:r>' + `
           '<w:proofErr w:type="spellStart"/><w:r><w:t>Genero</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
           '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
           '<w:proofErr w:type="spellStart"/><w:r><w:t>DataNascimento</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
           '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
           '<w:proofErr w:type="gramStart"/><w:r><w:t>pais-&gt;Pais</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
           '<w:r><w:t>)</w:t></w:r>'

$pkg1415 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p>' + $inner14 + '</w:p><w:p>' + $inner15 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full1415.InsertXML($pkg1415)
